# Generate Report for handback: append a new row for
# 30609edd-e53e-473b-ae91-2c1e3c3027da to the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$guid = "30609edd-e53e-473b-ae91-2c1e3c3027da"
$mdName = "$guid.md"
$statusInSync = "Handed back: in sync with en-US"
$include = "Include"

# ---------------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = $mdName
$wsOverview.Range("B4").Value = $statusInSync
$wsOverview.Range("C4").Value = $statusInSync

$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/HEAD/e2e/$mdName", "", "", $mdName)

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------------
$zhXlf = "$guid.5dfb08a6c2484f892d856d85ca00a2cdb26e6d29.zh-cn.xlf"
$zhHandoffDt = "2016-02-15 03:34:06"
$zhHandbackDt = "2016-02-15 03:34:51"

$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A4").Value = $mdName
$wsZh.Range("B4").Value = $statusInSync
$wsZh.Range("C4").Value = $zhXlf
$wsZh.Range("D4").Value = $zhHandoffDt
$wsZh.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("E4").Value = $mdName
$wsZh.Range("F4").Value = $zhXlf
$wsZh.Range("G4").Value = $zhHandbackDt
$wsZh.Range("H4").Value = $include

$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/HEAD/e2e/$mdName", "", "", $mdName)
$wsZh.Hyperlinks.Add($wsZh.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/HEAD/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zhXlf", "", "", $zhXlf)
$wsZh.Hyperlinks.Add($wsZh.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/HEAD/e2e/$mdName", "", "", $mdName)
$wsZh.Hyperlinks.Add($wsZh.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/HEAD/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zhXlf", "", "", $zhXlf)

# ---------------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------------
$deXlf = "$guid.5dfb08a6c2484f892d856d85ca00a2cdb26e6d29.de-de.xlf"
$deHandoffDt = "2016-02-15 03:34:19"
$deHandbackDt = "2016-02-15 03:35:17"

$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A4").Value = $mdName
$wsDe.Range("B4").Value = $statusInSync
$wsDe.Range("C4").Value = $deXlf
$wsDe.Range("D4").Value = $deHandoffDt
$wsDe.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("E4").Value = $mdName
$wsDe.Range("F4").Value = $deXlf
$wsDe.Range("G4").Value = $deHandbackDt
$wsDe.Range("H4").Value = $include

$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/HEAD/e2e/$mdName", "", "", $mdName)
$wsDe.Hyperlinks.Add($wsDe.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/HEAD/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$deXlf", "", "", $deXlf)
$wsDe.Hyperlinks.Add($wsDe.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/HEAD/e2e/$mdName", "", "", $mdName)
$wsDe.Hyperlinks.Add($wsDe.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/HEAD/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$deXlf", "", "", $deXlf)

Write-Output "done"
